$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": bump template version 1.0.0 -> 1.0.1 ---
$ws1 = $wb.Worksheets.Item("isa_template")
$ws1.Range("B4").Value = "1.0.1"

# --- Sheet "New Table": add example values to the mandatory ENA template row ---
$ws2 = $wb.Worksheets.Item("New Table")

$ws2.Range("B2").Value = "flowering stage"
$ws2.Range("C2").Value = "PO"
$ws2.Range("D2").Value = "http://purl.obolibrary.org/obo/PO_0007616"

# "2022-08-23" must stay a literal text string, not become a date serial.
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "2022-08-23"

$ws2.Range("H2").Value = "Germany"
$ws2.Range("I2").Value = "NCIT"
$ws2.Range("J2").Value = "http://purl.obolibrary.org/obo/NCIT_C16636"

# "+50.55" / "+6.21" must stay literal text, not become numbers (losing the "+").
$ws2.Range("K2").NumberFormat = "@"
$ws2.Range("K2").Value = "+50.55"
$ws2.Range("N2").NumberFormat = "@"
$ws2.Range("N2").Value = "+6.21"

$ws2.Range("Q2").Value = "petiole epidermis"
$ws2.Range("R2").Value = "PO"
$ws2.Range("S2").Value = "http://purl.obolibrary.org/obo/PO_0000051"
$ws2.Range("T2").Value = "soil"
$ws2.Range("U2").Value = "ENVO"
$ws2.Range("V2").Value = "http://purl.obolibrary.org/obo/ENVO_00001998"
$ws2.Range("W2").Value = "https://doi.org/10.1038/nature22971"
